# Add two new days of COVID-19 case/death data (July 4 & July 5, 2020) to the
# "Data" sheet, backfill the missing DHEC-sum figures for July 3 (row 121),
# extend the three shared-formula columns (D, F, J) to cover the new rows,
# update the chart's source ranges on the "Chart" sheet, and move the
# selection to the new last cell - matching the upstream "Data files updated
# July 5, 2020" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# Step 1: copy existing number formats onto the cells that are about to
# receive new data, so the new cells inherit the same styles as their
# neighbours (date format in column A, computed-ratio format in column J,
# plain numbers elsewhere).
# ---------------------------------------------------------------------------

# Row 121 gains E121/F121 (every other column in that row is already filled).
$ws.Range("E120").Copy()
$ws.Range("E121").PasteSpecial(-4122)
$ws.Range("F120").Copy()
$ws.Range("F121").PasteSpecial(-4122)

# Row 122 is a brand-new row with data in every used column.
$newRowCols = @("A", "B", "C", "D", "E", "F", "H", "I", "J")
foreach ($col in $newRowCols) {
    $ws.Range($col + "121").Copy()
    $ws.Range($col + "122").PasteSpecial(-4122)
}

# Row 123 is a brand-new row, but (like row 121 before this edit) it has no
# DHEC daily-case figure yet, so columns E/F stay blank.
$partialRowCols = @("A", "B", "C", "D", "H", "I", "J")
foreach ($col in $partialRowCols) {
    $ws.Range($col + "121").Copy()
    $ws.Range($col + "123").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: fill in the actual values / formulas.
# ---------------------------------------------------------------------------

# Row 121 - newly reported DHEC daily-case number and running sum.
$ws.Range("E121").Value = 1836
$ws.Range("F121").Formula = "=SUM(F120,E121)"

# Row 122 - July 4, 2020.
$ws.Range("A122").Value = 44016
$ws.Range("B122").Value = 43260
$ws.Range("C122").Value = 806
$ws.Range("D122").Formula = "=SUM(C122,-(C121))"
$ws.Range("E122").Value = 1463
$ws.Range("F122").Formula = "=SUM(F121,E122)"
$ws.Range("H122").Value = 56764
$ws.Range("I122").Value = 463201
$ws.Range("J122").Formula = "=IMDIV(H122,I122)"

# Row 123 - July 5, 2020.
$ws.Range("A123").Value = 44017
$ws.Range("B123").Value = 44717
$ws.Range("C123").Value = 813
$ws.Range("D123").Formula = "=SUM(C123,-(C122))"
$ws.Range("H123").Value = 58491
$ws.Range("I123").Value = 473543
$ws.Range("J123").Formula = "=IMDIV(H123,I123)"

# ---------------------------------------------------------------------------
# Step 3: move the active selection to the new last data cell.
# ---------------------------------------------------------------------------
$ws.Range("J123").Select()

# ---------------------------------------------------------------------------
# Step 4: point the chart's series at the widened Data range.
# ---------------------------------------------------------------------------
$chartWs = $wb.Worksheets.Item("Chart")
$chartObj = $chartWs.ChartObjects().Item(1)
$chart = $chartObj.Chart

$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = "=SERIES(Data!`$B`$1,Data!`$A`$2:`$A`$123,Data!`$B`$2:`$B`$123,1)"

$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = "=SERIES(Data!`$C`$1,Data!`$A`$2:`$A`$123,Data!`$C`$2:`$C`$123,2)"
